# Replace the ESPN-comparison table (NBA, Monday 12th Feb 2024 -> NBA, Tuesday 13th Feb 2024)
# New slate has 6 games (rows 2-7) instead of the old 10 games (rows 2-11).
# Values are written column-by-column (A2:A7, then B2:B7, then C2:C7, then the
# A1 header) to match the order new strings were originally authored in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: matchups
$ws.Range("A2").Value = "Boston Celtics (41-12) vs Brooklyn Nets (21-31)"
$ws.Range("A3").Value = "Oklahoma City Thunder (36-17) vs Orlando Magic (29-24)"
$ws.Range("A4").Value = "Miami Heat (28-25) vs Milwaukee Bucks (35-19)"
$ws.Range("A5").Value = "Sacramento Kings (30-22) vs Phoenix Suns (31-22)"
$ws.Range("A6").Value = "Minnesota Timberwolves (37-16) vs Portland Trail Blazers (15-37)"
$ws.Range("A7").Value = "Detroit Pistons (8-44) vs Los Angeles Lakers (29-26)"

# Column B: Ballgorithm picks
$ws.Range("B2").Value = "Boston Celtics (89.29%)"
$ws.Range("B3").Value = "Oklahoma City Thunder (77.78%)"
$ws.Range("B4").Value = "Milwaukee Bucks (79.31%)"
$ws.Range("B5").Value = "Sacramento Kings (62.50%)"
$ws.Range("B6").Value = "Minnesota Timberwolves (79.17%)"
$ws.Range("B7").Value = "Los Angeles Lakers (67.86%)"

# Column C: ESPN picks
$ws.Range("C2").Value = "Boston Celtics (71.9%)"
$ws.Range("C3").Value = "Oklahoma City Thunder (55.7%)"
$ws.Range("C4").Value = "Milwaukee Bucks (74.6%)"
$ws.Range("C5").Value = "Phoenix Suns (71.6%)"
$ws.Range("C6").Value = "Minnesota Timberwolves (81.8%)"
$ws.Range("C7").Value = "Los Angeles Lakers (76.7%)"

# Header (A1) with the new date, written last so the shared-string table
# lines up with how the workbook was originally regenerated.
$ws.Range("A1").Value = "NBA, Tuesday 13th Feb 2024"

# Drop the now-unused trailing rows (old table had 10 games / 11 rows total).
$ws.Rows("8:11").Delete()

# Reset the cursor back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
